$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data to append: column A = sequential index (204..215), column B = value
$newData = @(
    @(204, "-3.386324785396392E-17"),
    @(205, "-6.134613019590709E-17"),
    @(206, "-6.756747938929664E-17"),
    @(207, "-2.987579319737834E-17"),
    @(208, "5.063224145507306E-17"),
    @(209, "1.15421065562314E-16"),
    @(210, "-4.492211001331606E-17"),
    @(211, "2.687737185591566E-17"),
    @(212, "-3.361026734705064E-17"),
    @(213, "-1.345856296778673E-16"),
    @(214, "-6.732895491134983E-17"),
    @(215, "0")
)

$startRow = 206
$endRow = $startRow + $newData.Count - 1

# Replicate the direct cell formatting (bold/centered/bordered style) used by
# the existing column A cells onto the newly appended rows.
$ws.Range("A205").Copy()
for ($r = $startRow; $r -le $endRow; $r++) {
    $ws.Range("A$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $aVal = $newData[$i][0]
    $bVal = [double]$newData[$i][1]

    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 2).Value = $bVal
}
